$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 113, shifting the existing rows 113:155 down to 114:156
$ws.Rows.Item(113).Insert()

# Populate the newly inserted row with the new weekly record
$ws.Cells.Item(113, 1).Value = 9
$ws.Cells.Item(113, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(113, 3).Value = "Metropolitana"
$ws.Cells.Item(113, 4).Value = 44468
$ws.Cells.Item(113, 5).Value = 13
$ws.Cells.Item(113, 6).Value = 300000001
$ws.Cells.Item(113, 7).Value = "Rabanito"
$ws.Cells.Item(113, 8).Value = "Sin especificar"
$ws.Cells.Item(113, 9).Value = "Primera"
$ws.Cells.Item(113, 10).Value = 6100
$ws.Cells.Item(113, 11).Value = 3500
$ws.Cells.Item(113, 12).Value = 4000
$ws.Cells.Item(113, 13).Value = 3750
$ws.Cells.Item(113, 14).Value = "`$/cien unidades (volumen en unidades)"
$ws.Cells.Item(113, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(113, 16).Value = 38
$ws.Cells.Item(113, 17).Value = 100
$ws.Cells.Item(113, 18).Value = "Hortaliza"

Write-Output "Inserted row 113 and populated new record"
